# Extract short Chinese labels into the business_class column (renamed to
# business_short), replacing the raw numeric/"None" codes with a readable
# label looked up per business_id. Also fixes two "重打牌" class_name values
# to "重打牌登机牌" for the reprint-boarding-pass rows, and fills the
# class_name for the departure-query rows that previously read "None".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: business_class -> business_short
$ws.Range("C1").Value = "business_short"

# Map of (business_id, old business_class code) -> new short label.
# Built from the exact (B,C) pairs observed in the sheet.
$map = @{
    "01010102|1"   = "登录按钮"
    "01010903|5"   = "用户注销按钮"
    "01020100|81"  = "个人中心页面加载"
    "040O0238|9"   = "重打登机牌按钮"
    "040P0100|31"  = "检索页面加载"
    "042E0100|2"   = "旅客检索页面加载"
    "042E0115|2"   = "旅客提取页签检索按钮"
    "042E0117|5"   = "旅客检索清空按钮"
    "042E0118|2"   = "序号输入框选中旅客"
    "072K0202|None" = "出港查询按钮"
    "nav|5"        = "关闭当前页面"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 159) { $lastRow = 159 }

for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    if ($null -eq $bVal -or $null -eq $cVal) { continue }
    $key = "$bVal|$cVal"
    if ($map.ContainsKey($key)) {
        $ws.Cells.Item($r, 3).Value = $map[$key]
    }
}

# class_name (column E) fix-ups tied to the same rows.
for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    if ($bVal -eq "072K0202" -and $eVal -eq "None") {
        $ws.Cells.Item($r, 5).Value = "查询航班"
    }
    if ($bVal -eq "040O0238" -and $eVal -eq "重打牌") {
        $ws.Cells.Item($r, 5).Value = "重打牌登机牌"
    }
}
